# Updates Sales section for Examplery Data Excel sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Detailed Profit & Loss")
$ws.Activate()

# SALES OF GOODS/SERVICES (row 19)
$ws.Range("C19").Value = 10000
$ws.Range("D19").Value = 15000
$ws.Range("E19").Value = 20000
$ws.Range("G19").Value = 15000
$ws.Range("H19").Value = 15000
$ws.Range("I19").Value = 30000
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 30000
$ws.Range("M19").Value = 45000

# COMMISSIONS/FEES/ETC. (row 20)
$ws.Range("C20").Value = 3500
$ws.Range("D20").Value = 4000
$ws.Range("E20").Value = 4500
$ws.Range("G20").Value = 3500
$ws.Range("H20").Value = 4000
$ws.Range("I20").Value = 5500
$ws.Range("K20").Value = 3500
$ws.Range("L20").Value = 5500
$ws.Range("M20").Value = 6500

# SALES / OTHER (row 21)
$ws.Range("C21").Value = 500
$ws.Range("D21").Value = 1500
$ws.Range("E21").Value = 2000
$ws.Range("G21").Value = 2500
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 1200
$ws.Range("K21").Value = 2000
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = 25000

# Update the view: keep the header rows (1:2) frozen, scroll the
# lower pane down so row 12 is the first visible row, then move the
# selection to M22.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A3").Select()
$win.FreezePanes = $true
$win.ScrollRow = 12
$win.ScrollColumn = 1
$ws.Range("M22").Select()

$excel.Calculate()
